$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = 131143934
$ws.Range("B9").Value = 79243
$ws.Range("D9").Value = 'NT'
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = 'Garnlav'
$ws.Range("G9").Value = 'Alectoria sarmentosa'
$ws.Range("H9").Value = '(Ach.) Ach.'
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = ""
$ws.Range("I9").Style = "Normal"
$ws.Range("P9").Value = 'Storbackmyran, Mpd'
$ws.Range("Q9").Value = 562502
$ws.Range("R9").Value = 6917029
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = 'Västernorrland'
$ws.Range("U9").Value = 'Ånge'
$ws.Range("V9").Value = 'Medelpad'
$ws.Range("W9").Value = 'Torp'
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = '2026-01-03'
$ws.Range("Y9").Style = "Normal"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = '2026-01-03'
$ws.Range("AA9").Style = "Normal"
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AT9").NumberFormat = "@"
$ws.Range("AT9").Value = ""
$ws.Range("AT9").Style = "Normal"
$ws.Range("AW9").Value = 'Markus Borja'
$ws.Range("AX9").Value = 'Markus Borja'
$ws.Range("AY9").NumberFormat = "@"
$ws.Range("AY9").Value = ""
$ws.Range("AY9").Style = "Normal"

# Row 10
$ws.Range("A10").Value = 131143933
$ws.Range("B10").Value = 79243
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = 'Garnlav'
$ws.Range("G10").Value = 'Alectoria sarmentosa'
$ws.Range("H10").Value = '(Ach.) Ach.'
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = ""
$ws.Range("I10").Style = "Normal"
$ws.Range("P10").Value = 'Storbackmyran, Mpd'
$ws.Range("Q10").Value = 562495
$ws.Range("R10").Value = 6917033
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = 'Västernorrland'
$ws.Range("U10").Value = 'Ånge'
$ws.Range("V10").Value = 'Medelpad'
$ws.Range("W10").Value = 'Torp'
$ws.Range("Y10").NumberFormat = "@"
$ws.Range("Y10").Value = '2026-01-03'
$ws.Range("Y10").Style = "Normal"
$ws.Range("AA10").NumberFormat = "@"
$ws.Range("AA10").Value = '2026-01-03'
$ws.Range("AA10").Style = "Normal"
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AT10").NumberFormat = "@"
$ws.Range("AT10").Value = ""
$ws.Range("AT10").Style = "Normal"
$ws.Range("AW10").Value = 'Markus Borja'
$ws.Range("AX10").Value = 'Markus Borja'
$ws.Range("AY10").NumberFormat = "@"
$ws.Range("AY10").Value = ""
$ws.Range("AY10").Style = "Normal"

# Row 11
$ws.Range("A11").Value = 131143935
$ws.Range("B11").Value = 79243
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = 'Garnlav'
$ws.Range("G11").Value = 'Alectoria sarmentosa'
$ws.Range("H11").Value = '(Ach.) Ach.'
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = ""
$ws.Range("I11").Style = "Normal"
$ws.Range("P11").Value = 'Storbackmyran, Mpd'
$ws.Range("Q11").Value = 562507
$ws.Range("R11").Value = 6917024
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = 'Västernorrland'
$ws.Range("U11").Value = 'Ånge'
$ws.Range("V11").Value = 'Medelpad'
$ws.Range("W11").Value = 'Torp'
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("Y11").Value = '2026-01-03'
$ws.Range("Y11").Style = "Normal"
$ws.Range("AA11").NumberFormat = "@"
$ws.Range("AA11").Value = '2026-01-03'
$ws.Range("AA11").Style = "Normal"
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AG11").Value = $false
$ws.Range("AT11").NumberFormat = "@"
$ws.Range("AT11").Value = ""
$ws.Range("AT11").Style = "Normal"
$ws.Range("AW11").Value = 'Markus Borja'
$ws.Range("AX11").Value = 'Markus Borja'
$ws.Range("AY11").NumberFormat = "@"
$ws.Range("AY11").Value = ""
$ws.Range("AY11").Style = "Normal"
